$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.570.48"
$ws.Range("E2").Value = "  +1.89%  "
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.71"
$ws.Range("E5").Value = "  +1.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.22"
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.488.29"
$ws.Range("E8").Value = "  +1.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.590"
$ws.Range("E9").Value = "  +5.35%  "
$ws.Range("E10").Value = "  +4.58%  "
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").Value = "4.095.83"
$ws.Range("E13").Value = "  +1.46%  "
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.04"
$ws.Range("E15").Value = "  +2.29%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "66.604.18"
$ws.Range("E16").Value = "  +2.08%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000178"
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("D18").Value = "3.495.95"
$ws.Range("E18").Value = "  +1.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.27"
$ws.Range("E19").Value = "  +0.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.97"
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "391.50"
$ws.Range("E21").Value = "  +2.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.89"
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.79"
$ws.Range("E23").Value = "  +1.68%  "
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.531"
$ws.Range("E25").Value = "  +2.42%  "
$ws.Range("E26").Value = "  +1.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.16"
$ws.Range("E27").Value = "  +3.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.180"
$ws.Range("E28").Value = "  -0.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.33"
$ws.Range("E30").Value = "  +1.88%  "
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("E32").Value = "  +1.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.62"
$ws.Range("E33").Value = "  +1.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.30"
$ws.Range("E34").Value = "  +0.67%  "
$ws.Range("E35").Value = "  +3.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.42"
$ws.Range("E36").Value = "  +1.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.893"
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("E38").Value = "  +2.01%  "
$ws.Range("E39").Value = "  +2.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.63"
$ws.Range("E40").Value = "  +3.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0737"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.24"
$ws.Range("E42").Value = "  +0.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.72"
$ws.Range("E43").Value = "  +0.72%  "
$ws.Range("D44").Value = "2.773.66"
$ws.Range("E44").Value = "  -0.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.85"
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.54"
$ws.Range("E46").Value = "  +1.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0309"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "342.90"
$ws.Range("E48").Value = "  +2.43%  "
$ws.Range("E49").Value = "  +1.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.51"
$ws.Range("E50").Value = "  +6.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.853"
$ws.Range("E51").Value = "  +2.88%  "
